$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.883.96'
$ws.Range("E2").Value = '  +2.57%  '

$ws.Range("D3").Value = '2.112.49'
$ws.Range("E3").Value = '  +9.92%  '

$ws.Range("D4").Value = '''1.004'

$ws.Range("D5").Value = '''335.59'
$ws.Range("E5").Value = '  +5.02%  '

$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").Value = '''0.5297'
$ws.Range("E7").Value = '  +4.24%  '

$ws.Range("D8").Value = '''0.4355'
$ws.Range("E8").Value = '  +7.98%  '

$ws.Range("D9").Value = '''0.09010'
$ws.Range("E9").Value = '  +8.17%  '

$ws.Range("D10").Value = '''45.91'
$ws.Range("E10").Value = '  +8.77%  '

$ws.Range("D11").Value = '''1.176'
$ws.Range("E11").Value = '  +5.27%  '

$ws.Range("D12").Value = '''25.00'
$ws.Range("E12").Value = '  +3.53%  '

$ws.Range("D13").Value = '2.114.96'
$ws.Range("E13").Value = '  +10.29%  '

$ws.Range("D14").Value = '''6.747'
$ws.Range("E14").Value = '  +5.04%  '

$ws.Range("D15").Value = '''7.787'
$ws.Range("E15").Value = '  +7.34%  '

$ws.Range("D16").Value = '''97.47'
$ws.Range("E16").Value = '  +5.29%  '

$ws.Range("E17").Value = '  +0.22%  '

$ws.Range("D18").Value = '''0.00001133'
$ws.Range("E18").Value = '  +3.41%  '

$ws.Range("D19").Value = '''0.06665'
$ws.Range("E19").Value = '  +2.54%  '

$ws.Range("D20").Value = '''19.10'
$ws.Range("E20").Value = '  +3.28%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").Value = '''6.356'
$ws.Range("E22").Value = '  +6.81%  '

$ws.Range("D23").Value = '30.959.45'
$ws.Range("E23").Value = '  +2.78%  '

$ws.Range("E24").Value = '  +7.21%  '

$ws.Range("D25").Value = '2.361.58'
$ws.Range("E25").Value = '  +10.48%  '

$ws.Range("D26").Value = '''2.273'
$ws.Range("E26").Value = '  +3.66%  '

$ws.Range("D27").Value = '''22.77'
$ws.Range("E27").Value = '  +4.05%  '

$ws.Range("D28").Value = '''2.563'
$ws.Range("E28").Value = '  +12.62%  '

$ws.Range("D29").Value = '''163.05'
$ws.Range("E29").Value = '  +0.32%  '

$ws.Range("D30").Value = '''133.46'
$ws.Range("E30").Value = '  +3.33%  '

$ws.Range("D31").Value = '''1.167'
$ws.Range("E31").Value = '  +2.54%  '

$ws.Range("D32").Value = '''0.1072'
$ws.Range("E32").Value = '  +2.43%  '

$ws.Range("D33").Value = '''6.223'
$ws.Range("E33").Value = '  +4.49%  '

$ws.Range("D34").Value = '''4.015'

$ws.Range("D35").Value = '''1.520'
$ws.Range("E35").Value = '  +22.87%  '

$ws.Range("D36").Value = '''0.02616'
$ws.Range("E36").Value = '  +6.67%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''5.535'
$ws.Range("E37").Value = '  +4.15%  '

$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '''12.84'
$ws.Range("E38").Value = '  +11.32%  '

$ws.Range("D39").Value = '''9.512'
$ws.Range("E39").Value = '  +9.86%  '

$ws.Range("D40").Value = '''0.06718'
$ws.Range("E40").Value = '  +3.92%  '

$ws.Range("D41").Value = '''0.2267'
$ws.Range("E41").Value = '  +5.52%  '

$ws.Range("D42").Value = '''0.6857'
$ws.Range("E42").Value = '  +5.98%  '

$ws.Range("D43").Value = '''1.251'
$ws.Range("E43").Value = '  +2.99%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '''1.003'
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").Value = '''14.12'
$ws.Range("E45").Value = '  +6.17%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.6437'
$ws.Range("E46").Value = '  +6.37%  '

$ws.Range("D47").Value = '''2.240'
$ws.Range("E47").Value = '  +3.23%  '

$ws.Range("D48").Value = '''3.692'
$ws.Range("E48").Value = '  +1.94%  '

$ws.Range("D49").Value = '''1.275'
$ws.Range("E49").Value = '  +5.50%  '

$ws.Range("D50").Value = '''82.39'
$ws.Range("E50").Value = '  +5.65%  '

$ws.Range("D51").Value = '''1.166'
$ws.Range("E51").Value = '  +2.75%  '
